$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a weekly time series (rows 2..70) for Mango prices at
# "Vega Monumental Concepción". A new weekly observation (dated 44477, i.e.
# 2021-10-08) is inserted as a new record right before the old row 38
# (dated 44245), pushing that record and every one after it down by one row.
$ws.Rows.Item(38).Insert()

$ws.Cells.Item(38,1).Value = 11
$ws.Cells.Item(38,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(38,3).Value = "Bíobío"
$ws.Cells.Item(38,4).Value = 44477
$ws.Cells.Item(38,5).Value = 8
$ws.Cells.Item(38,6).Value = "Fruta"
$ws.Cells.Item(38,7).Value = 100108
$ws.Cells.Item(38,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(38,9).Value = 100108002
$ws.Cells.Item(38,10).Value = "Mango"
$ws.Cells.Item(38,11).Value = "Sin especificar"
$ws.Cells.Item(38,12).Value = "Primera"
$ws.Cells.Item(38,13).Value = 200
$ws.Cells.Item(38,14).Value = 8000
$ws.Cells.Item(38,15).Value = 8500
$ws.Cells.Item(38,16).Value = 8250
$ws.Cells.Item(38,17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(38,18).Value = "Perú"
$ws.Cells.Item(38,19).Value = 2062
$ws.Cells.Item(38,20).Value = 4
